$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 509.4717834624431
$ws.Range("C2").Value = 355.3164277330865
$ws.Range("D2").Value = 304.496326895058
$ws.Range("E2").Value = 280.0897893419756
$ws.Range("B3").Value = 598.9935381844241
$ws.Range("C3").Value = 419.091515081482
$ws.Range("D3").Value = 357.3594592318457
$ws.Range("E3").Value = 331.5630083491055
$ws.Range("B4").Value = 572.7150082520694
$ws.Range("C4").Value = 402.0101658226635
$ws.Range("D4").Value = 343.8982688741974
$ws.Range("E4").Value = 319.5015104924415
$ws.Range("B5").Value = 393.8699112426604
$ws.Range("C5").Value = 276.3110768592335
$ws.Range("D5").Value = 232.8408550402498
$ws.Range("E5").Value = 218.7053349206821
$ws.Range("B6").Value = 350.8963840745712
$ws.Range("C6").Value = 244.5030475833582
$ws.Range("D6").Value = 207.1892605313779
$ws.Range("E6").Value = 191.9988340830891
$ws.Range("B7").Value = 36.68411565055313
$ws.Range("C7").Value = 25.58486625401695
$ws.Range("D7").Value = 21.91596888786884
$ws.Range("E7").Value = 20.2530758127735
$ws.Range("B8").Value = 1998.468025715314
$ws.Range("C8").Value = 1393.8534122833
$ws.Range("D8").Value = 1204.279024121445
$ws.Range("E8").Value = 1103.372767904274
$ws.Range("B9").Value = 502.3898030613828
$ws.Range("C9").Value = 351.7173266962726
$ws.Range("D9").Value = 299.7340045080317
$ws.Range("E9").Value = 278.4070048960226
$ws.Range("B10").Value = 204.4680061353051
$ws.Range("C10").Value = 143.1438391229902
$ws.Range("D10").Value = 125.5911132001598
$ws.Range("E10").Value = 115.0712114082426
$ws.Range("B11").Value = 37.68742058254945
$ws.Range("C11").Value = 25.06104686309147
$ws.Range("D11").Value = 21.70492669859174
$ws.Range("E11").Value = 21.43795693297779
$ws.Range("B12").Value = 87.08378576801529
$ws.Range("C12").Value = 64.70116160461085
$ws.Range("D12").Value = 56.98604572718838
$ws.Range("E12").Value = 51.66515459616287
$ws.Range("B13").Value = 114.5930839945038
$ws.Range("C13").Value = 79.59768252295348
$ws.Range("D13").Value = 70.57154994444538
$ws.Range("E13").Value = 65.4454097161223
